$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Cells.Item(2, 1).Value = "Tổng công tại CẦN THƠ"
$ws.Cells.Item(2, 2).Value = 0

$ws.Cells.Item(3, 1).Value = "Phụ cấp tại CẦN THƠ"
$ws.Cells.Item(3, 2).Value = 0

$ws.Cells.Item(4, 1).Value = "Lương công tác tại CẦN THƠ"
$ws.Cells.Item(4, 2).Value = 0

$ws.Cells.Item(5, 1).Value = "Lương cơ bản tại CẦN THƠ"
$ws.Cells.Item(5, 2).Value = ""

$ws.Cells.Item(6, 1).Value = "Chiết khấu sale chính tại CẦN THƠ"
$ws.Cells.Item(6, 2).Value = 0

$ws.Cells.Item(7, 1).Value = "Chiết khấu sale phụ tại CẦN THƠ"
$ws.Cells.Item(7, 2).Value = 0

$ws.Cells.Item(8, 1).Value = "Đơn 1 bác sĩ tại CẦN THƠ"
$ws.Cells.Item(8, 2).Value = 0

$ws.Cells.Item(9, 1).Value = "Đơn 2 bác sĩ tại CẦN THƠ"
$ws.Cells.Item(9, 2).Value = 0

$ws.Cells.Item(10, 1).Value = "Công phụ phẫu 1 tại CẦN THƠ"
$ws.Cells.Item(10, 2).Value = 0

$ws.Cells.Item(11, 1).Value = "Công phụ phẫu 2 tại CẦN THƠ"
$ws.Cells.Item(11, 2).Value = 0

$ws.Cells.Item(12, 1).Value = "Chiết khấu thu nợ tại CẦN THƠ"
$ws.Cells.Item(12, 2).Value = 0

$ws.Cells.Item(13, 1).Value = "Ứng lương tại CẦN THƠ"
$ws.Cells.Item(13, 2).Value = -0

$ws.Cells.Item(14, 1).Value = "Tổng công tại LONG XUYÊN"
$ws.Cells.Item(14, 2).Value = 27.5

$ws.Cells.Item(15, 1).Value = "Phụ cấp tại LONG XUYÊN"
$ws.Cells.Item(15, 2).Value = 962500

$ws.Cells.Item(16, 1).Value = "Lương cơ bản tại LONG XUYÊN"
$ws.Cells.Item(16, 2).Value = 8328571.428571428

$ws.Cells.Item(17, 1).Value = "Chiết khấu sale chính tại LONG XUYÊN"
$ws.Cells.Item(17, 2).Value = 3035000

$ws.Cells.Item(18, 1).Value = "Chiết khấu sale phụ tại LONG XUYÊN"
$ws.Cells.Item(18, 2).Value = 0

$ws.Cells.Item(19, 1).Value = "Đơn 1 bác sĩ tại LONG XUYÊN"
$ws.Cells.Item(19, 2).Value = 0

$ws.Cells.Item(20, 1).Value = "Đơn 2 bác sĩ tại LONG XUYÊN"
$ws.Cells.Item(20, 2).Value = 0

$ws.Cells.Item(21, 1).Value = "Công phụ phẫu 1 tại LONG XUYÊN"
$ws.Cells.Item(21, 2).Value = 0

$ws.Cells.Item(22, 1).Value = "Công phụ phẫu 2 tại LONG XUYÊN"
$ws.Cells.Item(22, 2).Value = 0

$ws.Cells.Item(23, 1).Value = "Chiết khấu thu nợ tại LONG XUYÊN"
$ws.Cells.Item(23, 2).Value = 260000

$ws.Cells.Item(24, 1).Value = "Ứng lương tại LONG XUYÊN"
$ws.Cells.Item(24, 2).Value = -14000000

$ws.Cells.Item(25, 1).Value = "Tổng công tại SÓC TRĂNG"
$ws.Cells.Item(25, 2).Value = 0

$ws.Cells.Item(26, 1).Value = "Phụ cấp tại SÓC TRĂNG"
$ws.Cells.Item(26, 2).Value = 0

$ws.Cells.Item(27, 1).Value = "Lương công tác tại SÓC TRĂNG"
$ws.Cells.Item(27, 2).Value = 0

$ws.Cells.Item(28, 1).Value = "Lương cơ bản tại SÓC TRĂNG"
$ws.Cells.Item(28, 2).Value = ""

$ws.Cells.Item(29, 1).Value = "Chiết khấu sale chính tại SÓC TRĂNG"
$ws.Cells.Item(29, 2).Value = 0

$ws.Cells.Item(30, 1).Value = "Chiết khấu sale phụ tại SÓC TRĂNG"
$ws.Cells.Item(30, 2).Value = 0

$ws.Cells.Item(31, 1).Value = "Đơn 1 bác sĩ tại SÓC TRĂNG"
$ws.Cells.Item(31, 2).Value = 0

$ws.Cells.Item(32, 1).Value = "Đơn 2 bác sĩ tại SÓC TRĂNG"
$ws.Cells.Item(32, 2).Value = 0

$ws.Cells.Item(33, 1).Value = "Công phụ phẫu 1 tại SÓC TRĂNG"
$ws.Cells.Item(33, 2).Value = 0

$ws.Cells.Item(34, 1).Value = "Công phụ phẫu 2 tại SÓC TRĂNG"
$ws.Cells.Item(34, 2).Value = 0

$ws.Cells.Item(35, 1).Value = "Chiết khấu thu nợ tại SÓC TRĂNG"
$ws.Cells.Item(35, 2).Value = 0

$ws.Cells.Item(36, 1).Value = "Ứng lương tại SÓC TRĂNG"
$ws.Cells.Item(36, 2).Value = -0

$ws.Cells.Item(37, 1).Value = "Tổng lương tại CẦN THƠ"
$ws.Cells.Item(37, 2).Value = 0

$ws.Cells.Item(38, 1).Value = "Tổng lương tại LONG XUYÊN"
$ws.Cells.Item(38, 2).Value = -1413928.571428571

$ws.Cells.Item(39, 1).Value = "Tổng lương tại SÓC TRĂNG"
$ws.Cells.Item(39, 2).Value = 0

$ws.Cells.Item(40, 1).Value = "Tổng lương"
$ws.Cells.Item(40, 2).Value = -1413928.571428571
